$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# New engagement entries added for rows 18-22 (Date / Plik / Linie triplet in columns N/O/P)
$ws.Range("N18").Value = 45793
$ws.Range("O18").Value = "program.cs"
$ws.Range("P18").Value = 16

$ws.Range("N19").Value = 45793
$ws.Range("O19").Value = "appsettings.json"
$ws.Range("P19").Value = 7

$ws.Range("N20").Value = 45793
$ws.Range("O20").Value = "BadanieController.cs"
$ws.Range("P20").Value = 43

$ws.Range("N21").Value = 45793
$ws.Range("O21").Value = "LekarzController.cs"
$ws.Range("P21").Value = 37

$ws.Range("N22").Value = 45793
$ws.Range("O22").Value = "OsobaController.cs"
$ws.Range("P22").Value = 62

# Match the cell selection left active by the author when the file was saved
$ws.Range("M25").Select() | Out-Null

